$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new "price snapshot" column at AC, pushing the existing
# AC (nom) and AD (url_produit) columns one place to the right (-> AD, AE).
$ws.Columns("AC:AC").Insert()

# Header for the freshly inserted column.
$ws.Range("AC1").Value = "2026-01-28 22:16:04"

# Rows 2-80 already had a numeric price in column AB (the previous
# snapshot) - the new snapshot column repeats that same price.
for ($row = 2; $row -le 80; $row++) {
    $ws.Cells.Item($row, 29).Value = $ws.Cells.Item($row, 28).Value2
}

# Rows 81-205 had no price data (AB was blank) - the new snapshot column
# stays blank for them too.
for ($row = 81; $row -le 205; $row++) {
    $ws.Cells.Item($row, 29).Value = ""
}
